$d = $word.ActiveDocument

# --- "Versi" + "on"  ->  "Version" -------------------------------------
# A same-text assignment wouldn't register as an edit (no run merge), so
# briefly grow the range with an extra character spanning both runs, then
# remove that extra character again. Word folds the spanned runs into one
# as soon as a real text change touches both of them.
$rVersion = $d.Range(0, 7)
$rVersion.Text = "VersionX"
$d.Range(7, 8).Delete()

# --- " 2" + "."  ->  " 1."  (keeping the _GoBack bookmark in place) ----
# Growing the " 2" run by one character turns it into " 1." while leaving
# the bookmark and the old "." run completely untouched.
$rNum = $d.Range(7, 9)
$rNum.Text = " 1."

# The old "." run (now right after the bookmark) is redundant; remove it.
$d.Range(10, 11).Delete()
